# Regenerate orders with updated distance/size labels.
# The experiment's distance and size codes changed:
#   D80 -> D86, D51 -> D55, D64 -> D69, S30 -> S31
# These codes appear embedded inside many shared strings (condition
# labels, image filenames, the Distance column, and the Size column),
# so the safest way to reproduce the change through the Excel object
# model is a literal, case-sensitive substring Find & Replace over the
# whole used range for each old->new code pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# xlPart (2) = match substrings; xlByRows (1) = search order; MatchCase = $true
# so we only touch the exact tokens (the codes are already unambiguous,
# but being explicit keeps this safe if other text is added later).
$used.Replace("D80", "D86", 2, 1, $true) | Out-Null
$used.Replace("D51", "D55", 2, 1, $true) | Out-Null
$used.Replace("D64", "D69", 2, 1, $true) | Out-Null
$used.Replace("S30", "S31", 2, 1, $true) | Out-Null
